# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 11:45"

# --- Row 18: Banglades ---
$ws.Range("B18").Value = 352178
$ws.Range("C18").Value = 1557
$ws.Range("D18").Value = 260790
$ws.Range("E18").Value = 86381
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = 5007

# --- Row 25: Alemania ---
$ws.Range("B25").Value = 275581
$ws.Range("C25").Value = 30
$ws.Range("E25").Value = 19799
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 9482

# --- Row 47: Polonia ---
$ws.Range("B47").Value = 80699
$ws.Range("C47").Value = 711
$ws.Range("D47").Value = 64972
$ws.Range("E47").Value = 13411
$ws.Range("G47").Value = 18
$ws.Range("H47").Value = 2316

# --- Rows 98/99: Guinea & Malasia swap places (Malasia overtakes Guinea) ---
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 10358
$ws.Range("C98").Value = 82
$ws.Range("D98").Value = 9563
$ws.Range("E98").Value = 665
$ws.Range("H98").Value = 130

$ws.Range("A99").Value = "Guinea"
$ws.Range("B99").Value = 10344
$ws.Range("D99").Value = 9757
$ws.Range("E99").Value = 522
$ws.Range("H99").Value = 65

# --- Row 103: Finlandia ---
$ws.Range("B103").Value = 9195
$ws.Range("C103").Value = 149
$ws.Range("E103").Value = 1154

# --- Row 142: Sri Lanka ---
$ws.Range("D142").Value = 3118
$ws.Range("E142").Value = 168

# --- Row 144: Estonia ---
$ws.Range("B144").Value = 2976
$ws.Range("C144").Value = 36
$ws.Range("D144").Value = 2385
$ws.Range("E144").Value = 527

# --- Rows 214/215: Islas Malvinas & Montserrat swap places ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
